$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new user row (row 9): Niall / hilloe@mail.com / 1234
# Format C9 as text first so the numeric-looking password "1234" is
# stored as a string, matching the other password cells in the column.
$ws.Range("C9").NumberFormat = "@"

$ws.Range("A9").Value = "Niall"
$ws.Range("B9").Value = "hilloe@mail.com"
$ws.Range("C9").Value = "1234"
